$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain text values that often look
# numeric (e.g. "1.0000", "304.28", "23.345.93"). Excel auto-converts
# such literals to numbers when assigned directly, which would lose
# the original text formatting/precision. Temporarily force the
# column to Text format while writing, then clear the formatting
# override afterwards so the cells end up unstyled, matching the
# rest of the sheet.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '23.345.93'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '1.625.13'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '304.28'
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("D7").Value = '0.3780'
$ws.Range("D8").Value = '51.95'
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '1.227'
$ws.Range("E10").Value = '  -4.64%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.08092'
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("E13").Value = '  -2.64%  '
$ws.Range("D14").Value = '6.553'
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("E15").Value = '  -3.17%  '
$ws.Range("D16").Value = '7.221'
$ws.Range("E16").Value = '  -3.53%  '
$ws.Range("D17").Value = '1.624.69'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").Value = '93.48'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").Value = '0.06909'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("E20").Value = '  -3.28%  '
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("E22").Value = '  -2.83%  '
$ws.Range("D23").Value = '23.349.34'
$ws.Range("E23").Value = '  -0.90%  '
$ws.Range("D24").Value = '12.69'
$ws.Range("E24").Value = '  -2.54%  '
$ws.Range("D25").Value = '3.243'
$ws.Range("E25").Value = '  +3.56%  '
$ws.Range("D26").Value = '2.444'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("E27").Value = '  -1.84%  '
$ws.Range("D28").Value = '149.91'
$ws.Range("E28").Value = '  -1.37%  '
$ws.Range("D29").Value = '5.282'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = '134.22'
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").Value = '2.304'
$ws.Range("E31").Value = '  -4.89%  '
$ws.Range("D32").Value = '1.803.75'
$ws.Range("E32").Value = '  -0.65%  '
$ws.Range("D33").Value = '6.796'
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("E34").Value = '  +5.31%  '
$ws.Range("D35").Value = '0.9521'
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("D37").Value = '0.2512'
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").Value = '0.08821'
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").Value = '6.083'
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("D40").Value = '0.07126'
$ws.Range("E40").Value = '  -4.93%  '
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").Value = '0.7047'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("D43").Value = '16.13'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").Value = '12.30'
$ws.Range("E44").Value = '  -3.28%  '
$ws.Range("D45").Value = '0.6448'
$ws.Range("E45").Value = '  -3.13%  '
$ws.Range("D46").Value = '2.315'
$ws.Range("E46").Value = '  -2.41%  '
$ws.Range("D47").Value = '0.9993'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").Value = '3.991'
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").Value = '0.07976'
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("D51").Value = '125.63'
$ws.Range("E51").Value = '  -4.97%  '

$priceRange.ClearFormats()
